$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "A7"
$ws.Name = "A7"

# Update the item code column (A3:A9) from "A4" to "A7"
$ws.Range("A3:A9").Value = "A7"

# Add the new row 10: item code "A7", new account code "1234",
# and the description "Exchange Revenue / New Item Code"
$ws.Range("A9:C9").Copy($ws.Range("A10:C10"))

$ws.Range("A10").Value = "A7"

$ws.Range("B10").Formula = '="1234"'
$ws.Range("B10").Copy()
$ws.Range("B10").PasteSpecial(-4163)

$ws.Range("C10").Value = "Exchange Revenue / New Item Code"

$ws.Range("A10").Select()
